$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 13014.667
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 13014.667
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 13014.667
$ws.Range("M26").Value = $null
$ws.Range("N26").Value = -13702.667

$ws.Range("H40").Value = 3566.1667
$ws.Range("I40").Value = 3132.3333
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 3132.3333
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -2957.3333
$ws.Range("N40").Value = -4350

$ws.Range("H43").Value = 9717
$ws.Range("J43").Value = 7146.25
$ws.Range("L43").Value = 7146.25
$ws.Range("N43").Value = -7284.25

$ws.Range("H54").Value = 5000
$ws.Range("I54").Value = 5000
$ws.Range("K54").Value = 5000
$ws.Range("M54").Value = -4514

$ws.Range("H74").Value = 2862.5
$ws.Range("I74").Value = 2862.5
$ws.Range("K74").Value = 2862.5
$ws.Range("M74").Value = -1926.5

$ws.Range("H77").Value = 2862.5
$ws.Range("I77").Value = 2862.5
$ws.Range("K77").Value = 14312.5
$ws.Range("M77").Value = -9632.5

$ws.Range("H100").Value = 4671.4287
$ws.Range("I100").Value = 5033.3335
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 5033.3335
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -4492.3335
$ws.Range("N100").Value = -3582

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = $null
$ws.Range("N135").Value = $null

$ws.Range("H138").Value = 4677
$ws.Range("I138").Value = 4498.75
$ws.Range("K138").Value = 13496.25
$ws.Range("M138").Value = -8356.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 410.22223
$ws.Range("I5").Value = 406.14285
$ws.Range("J5").Value = 424.5
$ws.Range("K5").Value = 406.14285
$ws.Range("L5").Value = 424.5
$ws.Range("M5").Value = -294.14285
$ws.Range("N5").Value = -648.5

$ws.Range("H32").Value = 5825.364
$ws.Range("I32").Value = 5044.619
$ws.Range("K32").Value = 5044.619
$ws.Range("M32").Value = -4757.619

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = $null

$ws.Range("H50").Value = 14698.667
$ws.Range("J50").Value = 44000
$ws.Range("L50").Value = 44000
$ws.Range("N50").Value = -45428

$ws.Range("H132").Value = 3799.75
$ws.Range("I132").Value = 3799.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11399.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8869.25
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 410.22223
$ws.Range("I4").Value = 406.14285
$ws.Range("J4").Value = 424.5
$ws.Range("K4").Value = 406.14285
$ws.Range("L4").Value = 424.5
$ws.Range("M4").Value = -291.14285
$ws.Range("N4").Value = -654.5

$ws.Range("H20").Value = 2504.3125
$ws.Range("J20").Value = 4829.8
$ws.Range("L20").Value = 4829.8
$ws.Range("N20").Value = -5323.8

$ws.Range("H56").Value = 10110
$ws.Range("J56").Value = 10110
$ws.Range("L56").Value = 10110
$ws.Range("N56").Value = -11588

$ws.Range("H86").Value = 7499.25
$ws.Range("I86").Value = 7499.25
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 7499.25
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -6376.25
$ws.Range("N86").Value = $null

$ws.Range("H89").Value = 7499.25
$ws.Range("I89").Value = 7499.25
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 37496.25
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -31880.25
$ws.Range("N89").Value = $null

$ws.Range("H141").Value = 100526.336
$ws.Range("J141").Value = 100526.336
$ws.Range("L141").Value = 100526.336
$ws.Range("N141").Value = -110886.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 849.5
$ws.Range("I3").Value = 849.5
$ws.Range("K3").Value = 849.5
$ws.Range("M3").Value = -736.5

$ws.Range("H44").Value = 2500
$ws.Range("J44").Value = 2500
$ws.Range("L44").Value = 2500
$ws.Range("N44").Value = -3384

$ws.Range("H56").Value = 49996.5
$ws.Range("J56").Value = 49996.5
$ws.Range("L56").Value = 49996.5
$ws.Range("N56").Value = -51686.5

$ws.Range("H105").Value = 2099
$ws.Range("I105").Value = 1497.5
$ws.Range("K105").Value = 1497.5
$ws.Range("M105").Value = 249.5

$ws.Range("H107").Value = 928
$ws.Range("J107").Value = 874.5
$ws.Range("L107").Value = 874.5
$ws.Range("N107").Value = -4714.5

$ws.Range("H122").Value = 1810.5
$ws.Range("I122").Value = 1810.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5431.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2981.5
$ws.Range("N122").Value = $null

$ws.Range("H132").Value = 7710.7744
$ws.Range("J132").Value = 14999.833
$ws.Range("L132").Value = 44999.499
$ws.Range("N132").Value = -50059.499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1372.1666
$ws.Range("J131").Value = 1897.5
$ws.Range("L131").Value = 5692.5
$ws.Range("N131").Value = -15772.5

$ws.Range("H134").Value = 2222.5
$ws.Range("I134").Value = 2222.5
$ws.Range("K134").Value = 6667.5
$ws.Range("M134").Value = -1597.5

$ws.Range("H139").Value = 1040.8
$ws.Range("I139").Value = 1040.8
$ws.Range("K139").Value = 3122.4
$ws.Range("M139").Value = 2017.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 3000
$ws.Range("I21").Value = 3000
$ws.Range("K21").Value = 3000
$ws.Range("M21").Value = -2827

$ws.Range("H30").Value = 3000
$ws.Range("I30").Value = 3000
$ws.Range("K30").Value = 3000
$ws.Range("M30").Value = -2895

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = $null

$ws.Range("H54").Value = 9000
$ws.Range("I54").Value = 9000
$ws.Range("K54").Value = 9000
$ws.Range("M54").Value = -8610

$ws.Range("H132").Value = 1923.5938
$ws.Range("I132").Value = 1897.25
$ws.Range("J132").Value = 2108
$ws.Range("K132").Value = 5691.75
$ws.Range("L132").Value = 6324
$ws.Range("M132").Value = -3161.75
$ws.Range("N132").Value = -11384

$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1267.5
$ws.Range("I22").Value = 1250
$ws.Range("K22").Value = 1250
$ws.Range("M22").Value = -955

$ws.Range("H27").Value = 1267.5
$ws.Range("I27").Value = 1250
$ws.Range("K27").Value = 1250
$ws.Range("M27").Value = -1143

$ws.Range("H46").Value = 1849.7037
$ws.Range("J46").Value = 2283.0667
$ws.Range("L46").Value = 2283.0667
$ws.Range("N46").Value = -2659.0667

$ws.Range("H53").Value = 54999
$ws.Range("I53").Value = 54999
$ws.Range("K53").Value = 54999
$ws.Range("M53").Value = -54481

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = $null

$ws.Range("H63").Value = 43950
$ws.Range("I63").Value = 43950
$ws.Range("K63").Value = 43950
$ws.Range("M63").Value = -43201

$ws.Range("H66").Value = 43950
$ws.Range("I66").Value = 43950
$ws.Range("K66").Value = 131850
$ws.Range("M66").Value = -128106

$ws.Range("H132").Value = 2557
$ws.Range("I132").Value = 2568.5
$ws.Range("K132").Value = 7705.5
$ws.Range("M132").Value = -5175.5

$ws.Range("H136").Value = 4500
$ws.Range("J136").Value = 6000
$ws.Range("L136").Value = 18000
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = $null

$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -36108

$ws.Range("H132").Value = 67987.39999999999
$ws.Range("I132").Value = 67987.39999999999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 203962.2
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -201432.2
$ws.Range("N132").Value = $null
